$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) stays text so values like "1.000" or "5.220" are not
# coerced into numbers and lose their formatting (matches the original file,
# which stores these as literal strings, not numeric General-formatted cells).
$ws.Columns("D").NumberFormat = "@"

$ws.Range("D2").Value = "29.166.71"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").Value = "1.847.89"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "0.7034"
$ws.Range("E5").Value = "  -4.34%  "
$ws.Range("D6").Value = "238.63"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.3052"
$ws.Range("E8").Value = "  -3.25%  "
$ws.Range("D9").Value = "0.07423"
$ws.Range("E9").Value = "  +3.72%  "
$ws.Range("D10").Value = "23.41"
$ws.Range("E10").Value = "  -5.11%  "
$ws.Range("D11").Value = "0.08136"
$ws.Range("E11").Value = "  -2.09%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "0.7274"
$ws.Range("E12").Value = "  -3.70%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.220"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.818.25"
$ws.Range("E14").Value = "  -3.98%  "
$ws.Range("D15").Value = "88.90"
$ws.Range("E15").Value = "  -3.98%  "
$ws.Range("D16").Value = "29.649.67"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "5.779"
$ws.Range("E17").Value = "  -5.82%  "
$ws.Range("D18").Value = "238.70"
$ws.Range("E18").Value = "  -3.91%  "
$ws.Range("D19").Value = "13.08"
$ws.Range("E19").Value = "  -3.37%  "
$ws.Range("D20").Value = "0.000007653"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.174.05"
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "7.615"
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("D25").Value = "9.009"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("D26").Value = "160.76"
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("D27").Value = "0.1452"
$ws.Range("E27").Value = "  -7.06%  "
$ws.Range("D28").Value = "18.11"
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("D29").Value = "1.969"
$ws.Range("E29").Value = "  -3.60%  "
$ws.Range("D30").Value = "1.400"
$ws.Range("E30").Value = "  -5.00%  "
$ws.Range("D31").Value = "4.527"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D33").Value = "3.999"
$ws.Range("E33").Value = "  -4.46%  "
$ws.Range("D34").Value = "0.05185"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("D35").Value = "1.185"
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("D36").Value = "1.038"
$ws.Range("E36").Value = "  +4.11%  "
$ws.Range("D37").Value = "0.7045"
$ws.Range("E37").Value = "  -8.25%  "
$ws.Range("D39").Value = "0.01871"
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("D41").Value = "0.9465"
$ws.Range("E41").Value = "  +7.53%  "
$ws.Range("D42").Value = "6.014"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").Value = "0.4306"
$ws.Range("E43").Value = "  -5.86%  "
$ws.Range("D44").Value = "1.065.57"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("D45").Value = "70.24"
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "103.05"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.009.59"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.744"
$ws.Range("E49").Value = "  -5.76%  "
$ws.Range("D50").Value = "7.038"
$ws.Range("E50").Value = "  -6.70%  "
$ws.Range("D51").Value = "9.114"
$ws.Range("E51").Value = "  -4.48%  "
